$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: update date and reshuffled card contents ---
$ws.Range("A8").Value = 43970

$v_B8 = @'
type: signin
width: 2
height: 1
h3: Raise / Sponsor Funds
p: Try out our new feature. Raise Funds for your next project or Join us in building a brighter world.
button.primary: Create a Ticket*goto("/createticket")
button.secondary: View Tickets*goto("/tickets")
svg: /icons/stars.svg
'@
$ws.Range("B8").Value = $v_B8

$v_C8 = @'
type: signin
width: 2
height: 1
h3: Timeline
p: World history at a single place on a map. Scroll, Research, Paint and Memorize. Do not ever fall short on what happened at which time. Check it out.
button.default: Check it*goto("/timeline")
svg: /icons/bgtimeline.svg
'@
$ws.Range("C8").Value = $v_C8

$v_D8 = @'
type: blog
width: 2
height: 1
ser: 136
'@
$ws.Range("D8").Value = $v_D8

$v_E8 = @'
type: course
width: 2
height: 2
ser: 3,4,5,6,7
'@
$ws.Range("E8").Value = $v_E8

$v_F8 = @'
type: meetup
width: 2
height: 1
h3: Meetup coming in
date: 2020,5,7,10,30,0,0
button.default: Speak*goto("https://forms.gle/dyydXFRSsKzeH4hZ6")
button.default: Attend*goto("https://youtu.be/vscn-HP932E")
button.default: Details*goto("https://www.meetup.com/techshek/events/270179438/")
'@
$ws.Range("F8").Value = $v_F8

$v_G8 = @'
type: featured_blog
width: 2
height: 1
h3: Rules of being a good desi
p: Some rules to follow if you want to lit Pakistan brighter. We here at zakatlists are bounded by these rules. 😀
date: 6 Apr 2020
author: <a href=https://justaashir.com target=_blank>Aashir</a>
'@
$ws.Range("G8").Value = $v_G8

$v_H8 = @'
type: blog
width: 2
height: 1
ser: 135
'@
$ws.Range("H8").Value = $v_H8

$v_I8 = @'
type: subscribe
width: 2
height: 1
h3: Subscribe to stay tuned to zakatlists
input: enter your email here
button.default: Submit
'@
$ws.Range("I8").Value = $v_I8

$v_J8 = @'
type: blog
width: 2
height: 1
ser: 134
'@
$ws.Range("J8").Value = $v_J8

$v_K8 = @'
type: signin
width: 2
height: 1
h3.w-half: Sign up to get unlimited access to the entire content of zakatlists
button.primary: Sign In*goto("/signin/home")
button.secondary: Sign Up for Rs 300 / Month*goto("/signup")
'@
$ws.Range("K8").Value = $v_K8

$v_L8 = @'
type: footer
width: 6
height: 1
p.small: Eat from their fruits, and give the due alms on the day of harvest. <br> - Al Quran 6:141
facebook: https://facebook.com/zakatlists
twitter: https://twitter.com/zakatlists
makerlog: https://getmakerlog.com/@punch__lines 
'@
$ws.Range("L8").Value = $v_L8

# --- New column L (12th column) width ---
$ws.Columns(12).ColumnWidth = 35

# --- Row 9: brand-new row ---
$ws.Range("A9").Value = 43971
$ws.Range("A9").NumberFormat = $ws.Range("A8").NumberFormat

$v_B9 = @'
type: signin
width: 2
height: 1
h3: Raise / Sponsor Funds
p: Try out our new feature. Raise Funds for your next project or Join us in building a brighter world.
button.primary: Create a Ticket*goto("/createticket")
button.secondary: View Tickets*goto("/tickets")
svg: /icons/stars.svg
'@
$ws.Range("B9").Value = $v_B9
$ws.Range("B9").WrapText = $true

$v_C9 = @'
type: blog
width: 2
height: 1
ser: 136
'@
$ws.Range("C9").Value = $v_C9
$ws.Range("C9").WrapText = $true

$v_D9 = @'
type: course
width: 2
height: 2
ser: 3,4,5,6,7
'@
$ws.Range("D9").Value = $v_D9
$ws.Range("D9").WrapText = $true

$v_E9 = @'
type: meetup
width: 2
height: 1
h3: Meetup coming in
date: 2020,5,7,10,30,0,0
button.default: Speak*goto("https://forms.gle/dyydXFRSsKzeH4hZ6")
button.default: Attend*goto("https://youtu.be/vscn-HP932E")
button.default: Details*goto("https://www.meetup.com/techshek/events/270179438/")
'@
$ws.Range("E9").Value = $v_E9
$ws.Range("E9").WrapText = $true

$v_F9 = @'
type: featured_blog
width: 2
height: 1
h3: Rules of being a good desi
p: Some rules to follow if you want to lit Pakistan brighter. We here at zakatlists are bounded by these rules. 😀
date: 6 Apr 2020
author: <a href=https://justaashir.com target=_blank>Aashir</a>
'@
$ws.Range("F9").Value = $v_F9
$ws.Range("F9").WrapText = $true

$v_G9 = @'
type: blog
width: 2
height: 1
ser: 135
'@
$ws.Range("G9").Value = $v_G9
$ws.Range("G9").WrapText = $true

$v_H9 = @'
type: subscribe
width: 2
height: 1
h3: Subscribe to stay tuned to zakatlists
input: enter your email here
button.default: Submit
'@
$ws.Range("H9").Value = $v_H9
$ws.Range("H9").WrapText = $true

$v_I9 = @'
type: blog
width: 2
height: 1
ser: 134
'@
$ws.Range("I9").Value = $v_I9
$ws.Range("I9").WrapText = $true

$v_J9 = @'
type: signin
width: 2
height: 1
h3.w-half: Sign up to get unlimited access to the entire content of zakatlists
button.primary: Sign In*goto("/signin/home")
button.secondary: Sign Up for Rs 300 / Month*goto("/signup")
'@
$ws.Range("J9").Value = $v_J9
$ws.Range("J9").WrapText = $true

$v_K9 = @'
type: footer
width: 6
height: 1
p.small: Eat from their fruits, and give the due alms on the day of harvest. <br> - Al Quran 6:141
facebook: https://facebook.com/zakatlists
twitter: https://twitter.com/zakatlists
makerlog: https://getmakerlog.com/@punch__lines 
'@
$ws.Range("K9").Value = $v_K9
$ws.Range("K9").WrapText = $true

$ws.Rows(9).RowHeight = 255

# --- Update the active selection shown in the saved view ---
$ws.Range("L8").Select()
